# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet (fund-level detail) right before the
#    "总计" (totals) sheet, formatted like the existing quarterly sheets.
# 2) Add a new top row to the "总计" sheet summarizing the 2022-Q1 data,
#    pushing the existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q1" worksheet with per-fund holdings
# ---------------------------------------------------------------------
$srcSheet   = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

$ws = $wb.Worksheets.Add($totalSheet)
$ws.Name = "2022-Q1"

# NOTE: inserting a sheet shifts everybody after it, and a worksheet
# reference obtained *before* the insert tracks the tab *position*, not
# the sheet identity - so "$totalSheet" now actually points at the new
# "2022-Q1" sheet. Re-resolve "总计" by name after the insert.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy the layout/formatting of an existing fund-detail sheet so fonts,
# borders and column styles match the rest of the workbook.
$srcSheet.Range("A1:H11").Copy($ws.Range("A1"))
# The new sheet only has 4 funds (vs. 10 on 2021-Q4) - drop the extra rows.
$ws.Range("A6:H11").Clear()

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$ws.Range("B2").Value = "'377010"
$ws.Range("C2").Value = "上投摩根阿尔法混合"
$ws.Range("D2").Value = "'13.21"
$ws.Range("E2").Value = "'89.24"
$ws.Range("F2").Value = "'2.93"
$ws.Range("G2").Value = "'0.3871"
$ws.Range("H2").Value = 7

$ws.Range("B3").Value = "'000457"
$ws.Range("C3").Value = "上投摩根核心成长"
$ws.Range("D3").Value = "'12.84"
$ws.Range("E3").Value = "'87.63"
$ws.Range("F3").Value = "'2.97"
$ws.Range("G3").Value = "'0.3813"
$ws.Range("H3").Value = 8

$ws.Range("B4").Value = "'005526"
$ws.Range("C4").Value = "工银瑞信新生代消费灵活配置混合"
$ws.Range("D4").Value = "'2.23"
$ws.Range("E4").Value = "'94.09"
$ws.Range("F4").Value = "'4.10"
$ws.Range("G4").Value = "'0.0914"
$ws.Range("H4").Value = 5

$ws.Range("B5").Value = "'002005"
$ws.Range("C5").Value = "工银瑞信新得利混合"
$ws.Range("D5").Value = "'3.82"
$ws.Range("E5").Value = "'29.93"
$ws.Range("F5").Value = "'1.41"
$ws.Range("G5").Value = "'0.0539"
$ws.Range("H5").Value = 9

# ---------------------------------------------------------------------
# Part 2: prepend a 2022-Q1 summary row onto the "总计" sheet, shifting
# the older quarters down by one row.
# ---------------------------------------------------------------------
$tot = $totalSheet

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 4
$tot.Range("D2").Value = 0.91

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2021-Q4"
$tot.Range("C3").Value = 10
$tot.Range("D3").Value = 2.65

$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2021-Q3"
$tot.Range("C4").Value = 9
$tot.Range("D4").Value = 1.62

$tot.Range("A5").Value = 3
$tot.Range("B5").Value = "2021-Q2"
$tot.Range("C5").Value = 13
$tot.Range("D5").Value = 1.47

$tot.Range("A6").Value = 4
$tot.Range("B6").Value = "2021-Q1"
$tot.Range("C6").Value = 4
$tot.Range("D6").Value = 0.39

# The new row 6 (A6) needs the same formatting as the rest of column A
# (bold, thin border, centered) - apply it explicitly since it's a
# brand-new cell outside the sheet's original used range.
$tot.Range("A6").Font.Bold = $true
$tot.Range("A6").HorizontalAlignment = -4108
$tot.Range("A6").VerticalAlignment = -4160
$tot.Range("A6").Borders.LineStyle = 1
